$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.71511833333333
$ws.Range("H2").Value = 77.145355
$ws.Range("I2").Value = 0.5736986116453374
$ws.Range("J2").Value = 0.5736986116453374
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.33599166666667
$ws.Range("N2").Value = 52.007975
$ws.Range("O2").Value = 0.4573561888773979
$ws.Range("P2").Value = 0.4573561888773979
$ws.Range("Q2").Value = 445.7970771340139
$ws.Range("R2").Value = 4012.173694206125
$ws.Range("S2").Value = 0.2623846105863659
$ws.Range("T2").Value = 0.2623846105863659

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.71511833333333
$ws.Range("H3").Value = 77.145355
$ws.Range("I3").Value = 0.5736986116453374
$ws.Range("J3").Value = 0.5736986116453374
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.077707333333334
$ws.Range("N3").Value = 27.233122
$ws.Range("O3").Value = 0.2394870573052156
$ws.Range("P3").Value = 0.2394870573052156
$ws.Range("Q3").Value = 233.4343182720345
$ws.Range("R3").Value = 2100.90886444831
$ws.Range("S3").Value = 0.1373933922830295
$ws.Range("T3").Value = 0.1373933922830295

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.71511833333333
$ws.Range("H4").Value = 77.145355
$ws.Range("I4").Value = 0.5736986116453374
$ws.Range("J4").Value = 0.5736986116453374
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.491094
$ws.Range("N4").Value = 34.473282
$ws.Range("O4").Value = 0.3031567538173866
$ws.Range("P4").Value = 0.3031567538173866
$ws.Range("Q4").Value = 295.4948419894566
$ws.Range("R4").Value = 2659.45357790511
$ws.Range("S4").Value = 0.173920608775942
$ws.Range("T4").Value = 0.173920608775942

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.915995333333335
$ws.Range("H5").Value = 29.747986
$ws.Range("I5").Value = 0.2212236662524262
$ws.Range("J5").Value = 0.2212236662524262
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.33599166666667
$ws.Range("N5").Value = 52.007975
$ws.Range("O5").Value = 0.4573561888773979
$ws.Range("P5").Value = 0.4573561888773979
$ws.Range("Q5").Value = 171.9036124653723
$ws.Range("R5").Value = 1547.13251218835
$ws.Range("S5").Value = 0.1011780128866951
$ws.Range("T5").Value = 0.1011780128866951

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.915995333333335
$ws.Range("H6").Value = 29.747986
$ws.Range("I6").Value = 0.2212236662524262
$ws.Range("J6").Value = 0.2212236662524262
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.077707333333334
$ws.Range("N6").Value = 27.233122
$ws.Range("O6").Value = 0.2394870573052156
$ws.Range("P6").Value = 0.2394870573052156
$ws.Range("Q6").Value = 90.01450355469913
$ws.Range("R6").Value = 810.1305319922922
$ws.Range("S6").Value = 0.05298020483706468
$ws.Range("T6").Value = 0.05298020483706468

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.915995333333335
$ws.Range("H7").Value = 29.747986
$ws.Range("I7").Value = 0.2212236662524262
$ws.Range("J7").Value = 0.2212236662524262
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.491094
$ws.Range("N7").Value = 34.473282
$ws.Range("O7").Value = 0.3031567538173866
$ws.Range("P7").Value = 0.3031567538173866
$ws.Range("Q7").Value = 113.9456344788947
$ws.Range("R7").Value = 1025.510710310052
$ws.Range("S7").Value = 0.06706544852866646
$ws.Range("T7").Value = 0.06706544852866647

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.192279333333332
$ws.Range("H8").Value = 27.576838
$ws.Range("I8").Value = 0.2050777221022365
$ws.Range("J8").Value = 0.2050777221022365
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.33599166666667
$ws.Range("N8").Value = 52.007975
$ws.Range("O8").Value = 0.4573561888773979
$ws.Range("P8").Value = 0.4573561888773979
$ws.Range("Q8").Value = 159.3572779203389
$ws.Range("R8").Value = 1434.21550128305
$ws.Range("S8").Value = 0.09379356540433699
$ws.Range("T8").Value = 0.09379356540433699

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.192279333333332
$ws.Range("H9").Value = 27.576838
$ws.Range("I9").Value = 0.2050777221022365
$ws.Range("J9").Value = 0.2050777221022365
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.077707333333334
$ws.Range("N9").Value = 27.233122
$ws.Range("O9").Value = 0.2394870573052156
$ws.Range("P9").Value = 0.2394870573052156
$ws.Range("Q9").Value = 83.44482151424845
$ws.Range("R9").Value = 751.0033936282359
$ws.Range("S9").Value = 0.04911346018512139
$ws.Range("T9").Value = 0.0491134601851214

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.192279333333332
$ws.Range("H10").Value = 27.576838
$ws.Range("I10").Value = 0.2050777221022365
$ws.Range("J10").Value = 0.2050777221022365
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.491094
$ws.Range("N10").Value = 34.473282
$ws.Range("O10").Value = 0.3031567538173866
$ws.Range("P10").Value = 0.3031567538173866
$ws.Range("Q10").Value = 105.6293458935906
$ws.Range("R10").Value = 950.6641130423158
$ws.Range("S10").Value = 0.06217069651277813
$ws.Range("T10").Value = 0.06217069651277814

